$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")

# Add a new worksheet right after Sheet2 -> becomes "Sheet3" and the active tab
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws2)

# Populate the new "Sheet3" with the reorganized take/reserve/buy summary table
$ws3.Range("A1").Value = "take"
$ws3.Range("B1").Value = 15
$ws3.Range("C1").Value = "5 take 1, 5 take 2, 5 discard"

$ws3.Range("A2").Value = "Reserve"
$ws3.Range("B2").Value = 12
$ws3.Range("C2").Value = "any of 12"

$ws3.Range("B3").Value = 3
$ws3.Range("C3").Value = "top of deck"

$ws3.Range("A4").Value = "buy"
$ws3.Range("B4").Value = 12
$ws3.Range("C4").Value = "any of 12"

$ws3.Range("B5").Value = 3
$ws3.Range("C5").Value = "reserves"

$ws3.Range("B6").Value = 15
$ws3.Range("C6").Value = "with gold"

$ws3.Range("B7").Formula = "=SUM(B1:B6)"

# Update the selection left behind on Sheet2 (no longer the active tab)
$ws2.Range("A7:C11").Select() | Out-Null

# Sheet3 is now the active sheet/tab; leave its selection on B8
$ws3.Range("B8").Select() | Out-Null
$ws3.Activate()
